# Update betting-odds values for the "Jogos da Semana" worksheet.
# The diff only changes numeric values (no structural/formatting changes),
# so we simply write the new values into the referenced cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "G2"  = 2.3
    "I2"  = 2.9
    "J2"  = 2.88
    "K2"  = 2.3
    "L2"  = 3.4
    "O2"  = 1.22
    "P2"  = 4
    "U2"  = 1.62
    "V2"  = 2.2
    "Z2"  = 21
    "AI2" = 17
    "AO2" = 12
    "AP2" = 19
    "AU2" = 7.5

    "G3"  = 1.38
    "H3"  = 4.75
    "J3"  = 1.91
    "N3"  = 13
    "O3"  = 1.25
    "P3"  = 3.75
    "Q3"  = 1.8
    "R3"  = 2
    "S3"  = 1.36
    "T3"  = 3
    "U3"  = 2.1
    "V3"  = 1.67
    "W3"  = 6.5
    "AD3" = 9
    "AE3" = 23
    "AG3" = 451
    "AH3" = 17
    "AJ3" = 21
    "AK3" = 81
    "AS3" = 151
    "AT3" = 3
    "AW3" = 8.5
    "BA3" = 201

    "N4"  = 13.8

    "G6"  = 1.91
    "H6"  = 3.3
    "I6"  = 4
    "J6"  = 2.63
    "K6"  = 2.1
    "L6"  = 4.5
    "M6"  = 1.07
    "N6"  = 9
    "Q6"  = 2.1
    "R6"  = 1.7
    "X6"  = 8.5
    "Z6"  = 17
    "AA6" = 17
    "AB6" = 29
    "AC6" = 8.5
    "AD6" = 6.5
    "AH6" = 10
    "AI6" = 21
    "AN6" = 3.75
    "AO6" = 11
    "AR6" = 51
    "AS6" = 151
    "AV6" = 51
    "AW6" = 6
    "AX6" = 23
    "AY6" = 34
    "AZ6" = 81
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
